$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used in the data rows (A..T)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Row block reordering: data rows 2-35 (the weekly price blocks) are being
# rearranged into a new row order. Build a map of new-row -> old-row.
$rowMap = @{
    2  = 30
    3  = 31
    4  = 32
    5  = 14
    6  = 15
    7  = 16
    8  = 27
    9  = 28
    10 = 29
    11 = 12
    12 = 13
    13 = 17
    14 = 18
    15 = 19
    16 = 5
    17 = 6
    18 = 7
    19 = 33
    20 = 34
    21 = 35
    22 = 24
    23 = 25
    24 = 26
    25 = 2
    26 = 3
    27 = 4
    28 = 20
    29 = 21
    30 = 22
    31 = 23
    32 = 8
    33 = 9
    34 = 10
    35 = 11
}

# 1) Snapshot every source row (2-35) fully before writing anything, since the
#    permutation has overlapping cycles and writing in place would clobber
#    values that are still needed as a source for a later row.
$snapshot = @{}
for ($r = 2; $r -le 35; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range($col + $r).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshot back out into the new row positions.
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $srcRowVals = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range($col + $newRow).Value = $srcRowVals[$col]
    }
}
